$wb = $excel.ActiveWorkbook

# This script applies updated market-price figures to the per-job "Leve Profit"
# tables (columns H:N) on each worksheet, as produced by the scheduled pricing
# runner. Values are written directly (source data has no formulas).

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 847.4
$ws.Range("J2").Value = 1404.75
$ws.Range("L2").Value = 1404.75
$ws.Range("N2").Value = -1630.75
$ws.Range("H38").Value = 431.44446
$ws.Range("I38").Value = 360.375
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 1081.125
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -709.125
$ws.Range("N38").Value = -3744
$ws.Range("H53").Value = 94.666664
$ws.Range("I53").Value = 93.59999999999999
$ws.Range("K53").Value = 93.59999999999999
$ws.Range("M53").Value = 543.4
$ws.Range("H62").Value = 6041.0527
$ws.Range("I62").Value = 5534.5713
$ws.Range("K62").Value = 5534.5713
$ws.Range("M62").Value = -4910.5713
$ws.Range("H65").Value = 6041.0527
$ws.Range("I65").Value = 5534.5713
$ws.Range("K65").Value = 27672.8565
$ws.Range("M65").Value = -24552.8565
$ws.Range("H74").Value = 3787.3076
$ws.Range("I74").Value = 3787.3076
$ws.Range("K74").Value = 3787.3076
$ws.Range("M74").Value = -2851.3076
$ws.Range("H77").Value = 3787.3076
$ws.Range("I77").Value = 3787.3076
$ws.Range("K77").Value = 18936.538
$ws.Range("M77").Value = -14256.538
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 500
$ws.Range("M107").Value = 1420
$ws.Range("H132").Value = 3108.9807
$ws.Range("I132").Value = 2749.4468
$ws.Range("K132").Value = 8248.340400000001
$ws.Range("M132").Value = -5718.340400000001
$ws.Range("H135").Value = 2467.1667
$ws.Range("I135").Value = 2728.7
$ws.Range("J135").Value = 1159.5
$ws.Range("K135").Value = 24558.3
$ws.Range("L135").Value = 10435.5
$ws.Range("M135").Value = -22023.3
$ws.Range("N135").Value = -15505.5
$ws.Range("H137").Value = 2077.3333
$ws.Range("I137").Value = 1639.2
$ws.Range("K137").Value = 4917.6
$ws.Range("M137").Value = -2367.6
$ws.Range("H138").Value = 2632.4167
$ws.Range("I138").Value = 1329.5625
$ws.Range("J138").Value = 3004.6606
$ws.Range("K138").Value = 3988.6875
$ws.Range("L138").Value = 9013.981800000001
$ws.Range("M138").Value = 1151.3125
$ws.Range("N138").Value = -19293.9818
$ws.Range("H141").Value = 2910.0715
$ws.Range("I141").Value = 2849.3076
$ws.Range("J141").Value = 3700
$ws.Range("K141").Value = 8547.9228
$ws.Range("L141").Value = 11100
$ws.Range("M141").Value = -3367.9228
$ws.Range("N141").Value = -21460

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 5121.4
$ws.Range("I19").Value = 199.5
$ws.Range("K19").Value = 199.5
$ws.Range("M19").Value = 29.5
$ws.Range("H25").Value = 348.33334
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 348.33334
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 348.33334
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -1152.33334
$ws.Range("H61").Value = 5148.645
$ws.Range("I61").Value = 4170.778
$ws.Range("J61").Value = 11749.25
$ws.Range("K61").Value = 4170.778
$ws.Range("L61").Value = 11749.25
$ws.Range("M61").Value = -3958.778
$ws.Range("N61").Value = -12173.25
$ws.Range("H74").Value = 4040.2
$ws.Range("I74").Value = 3727
$ws.Range("J74").Value = 4845.5713
$ws.Range("K74").Value = 3727
$ws.Range("L74").Value = 4845.5713
$ws.Range("M74").Value = -2853
$ws.Range("N74").Value = -6593.5713
$ws.Range("H77").Value = 4040.2
$ws.Range("I77").Value = 3727
$ws.Range("J77").Value = 4845.5713
$ws.Range("K77").Value = 18635
$ws.Range("L77").Value = 24227.8565
$ws.Range("M77").Value = -14267
$ws.Range("N77").Value = -32963.85649999999
$ws.Range("H122").Value = 4350.104
$ws.Range("I122").Value = 3814.9512
$ws.Range("K122").Value = 11444.8536
$ws.Range("M122").Value = -8994.8536
$ws.Range("H136").Value = 5148.645
$ws.Range("I136").Value = 4170.778
$ws.Range("J136").Value = 11749.25
$ws.Range("K136").Value = 12512.334
$ws.Range("L136").Value = 35247.75
$ws.Range("M136").Value = -9962.334000000001
$ws.Range("N136").Value = -40347.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 22885.2
$ws.Range("I75").Value = 22885.2
$ws.Range("K75").Value = 22885.2
$ws.Range("M75").Value = -21949.2
$ws.Range("H78").Value = 22885.2
$ws.Range("I78").Value = 22885.2
$ws.Range("K78").Value = 68655.60000000001
$ws.Range("M78").Value = -63975.60000000001
$ws.Range("H80").Value = 960.5
$ws.Range("I80").Value = 138
$ws.Range("K80").Value = 138
$ws.Range("M80").Value = 860
$ws.Range("H82").Value = 5178.5
$ws.Range("I82").Value = 5178.5
$ws.Range("K82").Value = 5178.5
$ws.Range("M82").Value = -4795.5
$ws.Range("H83").Value = 960.5
$ws.Range("I83").Value = 138
$ws.Range("K83").Value = 690
$ws.Range("M83").Value = 4302
$ws.Range("H85").Value = 5178.5
$ws.Range("I85").Value = 5178.5
$ws.Range("K85").Value = 5178.5
$ws.Range("M85").Value = -3852.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4456.875
$ws.Range("I58").Value = 3958.3333
$ws.Range("K58").Value = 3958.3333
$ws.Range("M58").Value = -3755.3333
$ws.Range("H99").Value = 5664.722
$ws.Range("I99").Value = 4225.364
$ws.Range("J99").Value = 7926.5713
$ws.Range("K99").Value = 4225.364
$ws.Range("L99").Value = 7926.5713
$ws.Range("M99").Value = -2727.364
$ws.Range("N99").Value = -10922.5713
$ws.Range("H126").Value = 5664.722
$ws.Range("I126").Value = 4225.364
$ws.Range("J126").Value = 7926.5713
$ws.Range("K126").Value = 12676.092
$ws.Range("L126").Value = 23779.7139
$ws.Range("M126").Value = -10206.092
$ws.Range("N126").Value = -28719.7139
$ws.Range("H131").Value = 87298.5
$ws.Range("J131").Value = 86442.78
$ws.Range("L131").Value = 86442.78
$ws.Range("N131").Value = -96522.78
$ws.Range("H132").Value = 1567.5454
$ws.Range("I132").Value = 1225.5
$ws.Range("K132").Value = 3676.5
$ws.Range("M132").Value = -1146.5
$ws.Range("H136").Value = 4456.875
$ws.Range("I136").Value = 3958.3333
$ws.Range("K136").Value = 11874.9999
$ws.Range("M136").Value = -9324.999899999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2970.8572
$ws.Range("J59").Value = 4249.75
$ws.Range("L59").Value = 12749.25
$ws.Range("N59").Value = -13829.25
$ws.Range("H92").Value = 735
$ws.Range("J92").Value = 735
$ws.Range("L92").Value = 2205
$ws.Range("N92").Value = -4701

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 10248199
$ws.Range("I14").Value = 16740333
$ws.Range("K14").Value = 16740333
$ws.Range("M14").Value = -16740165
$ws.Range("H18").Value = 4007501
$ws.Range("I18").Value = 505002.5
$ws.Range("K18").Value = 505002.5
$ws.Range("M18").Value = -504709.5
$ws.Range("H80").Value = 7241.4287
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 9338
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 9338
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -11334
$ws.Range("H83").Value = 7241.4287
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 9338
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 46690
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -56674
$ws.Range("H102").Value = 1869.3226
$ws.Range("I102").Value = 1421.1538
$ws.Range("J102").Value = 4199.8
$ws.Range("K102").Value = 1421.1538
$ws.Range("L102").Value = 4199.8
$ws.Range("M102").Value = 200.8462
$ws.Range("N102").Value = -7443.8
$ws.Range("H126").Value = 7876.696
$ws.Range("J126").Value = 9617.5
$ws.Range("L126").Value = 28852.5
$ws.Range("N126").Value = -33792.5
$ws.Range("H132").Value = 6146.4287
$ws.Range("I132").Value = 5273.077
$ws.Range("K132").Value = 15819.231
$ws.Range("M132").Value = -13289.231
$ws.Range("H136").Value = 30897.295
$ws.Range("J136").Value = 30897.295
$ws.Range("L136").Value = 92691.88499999999
$ws.Range("N136").Value = -97791.88499999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H40").Value = 6943.4443
$ws.Range("J40").Value = 5917.75
$ws.Range("L40").Value = 5917.75
$ws.Range("N40").Value = -6189.75
$ws.Range("H122").Value = 3081.6365
$ws.Range("I122").Value = 2362.25
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7086.75
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4636.75
$ws.Range("N122").Value = -19900

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3353.4546
$ws.Range("I126").Value = 1498.625
$ws.Range("J126").Value = 8299.666999999999
$ws.Range("K126").Value = 4495.875
$ws.Range("L126").Value = 24899.001
$ws.Range("M126").Value = -2025.875
$ws.Range("N126").Value = -29839.001
